$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "City" column (E), shifting
# City -> F and Fund -> G. The new column becomes "Primary Email".
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Header + values for the new "Primary Email" column
$ws.Range("E1").Value = "Primary Email"
$ws.Range("E2").Value = "emp1@gmail.com"
$ws.Range("E3").Value = "emp2@gmail.com"
$ws.Range("E4").Value = "emp3@gmail.com"
$ws.Range("E5").Value = "emp4@gmail.com"
$ws.Range("E6").Value = "emp5@gmail.com"
$ws.Range("E7").Value = "emp6@gmail.com"

# Normalize the (duplicated) style used by column B's data cells so it
# matches the other body-style cells, and drop the leftover formatted
# but empty cells below the table.
$ws.Range("B1:B7").Style = "Normal 2"
$ws.Range("B8:B9").Clear()

$ws.Range("E1:E7").Select() | Out-Null
